$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Raw": insert the new drain_vec_collect_map_filter.rs benchmark
# block (2 rows: input-len 65536 and 262144) right after the existing
# collect_map_filter.rs block (rows 8:9), i.e. before the former row 10
# (collect_map.rs). This pushes all following rows down by 2.
# ---------------------------------------------------------------------
$raw = $wb.Worksheets.Item("Raw")
$raw.Rows("10:11").Insert()

# Row 10 — input len 65536
$raw.Range("B10").Value2 = "drain_vec_collect_map_filter.rs"
$raw.Range("C10").Value2 = "``inputs.par_drain().map(map).filter(filter).collect()``"
$raw.Range("D10").Value2 = 65536
$raw.Range("E10").Value2 = 10.45
$raw.Range("F10").Value2 = 12.68
$raw.Range("G10").Value2 = 3.81
$raw.Range("H10").Value2 = 4.12
$raw.Range("I10").Formula = '=IF(E10="","",TEXT(E10, "0.00")&" (1.00)")'
$raw.Range("J10").Formula = '=IF(F10="","",TEXT(F10,"0.00")&" ("&TEXT(F10/$E10,"0.00")&")")'
$raw.Range("K10").Formula = '=IF(G10="","",TEXT(G10,"0.00")&" ("&TEXT(G10/$E10,"0.00")&")")'
$raw.Range("L10").Formula = '=IF(H10="","",TEXT(H10,"0.00")&" ("&TEXT(H10/$E10,"0.00")&")")'

# Row 11 — input len 262144
$raw.Range("D11").Value2 = 262144
$raw.Range("E11").Value2 = 58.37
$raw.Range("F11").Value2 = 23.1
$raw.Range("G11").Value2 = 13.96
$raw.Range("H11").Value2 = 12.84
$raw.Range("I11").Formula = '=IF(E11="","",TEXT(E11, "0.00")&" (1.00)")'
$raw.Range("J11").Formula = '=IF(F11="","",TEXT(F11,"0.00")&" ("&TEXT(F11/$E11,"0.00")&")")'
$raw.Range("K11").Formula = '=IF(G11="","",TEXT(G11,"0.00")&" ("&TEXT(G11/$E11,"0.00")&")")'
$raw.Range("L11").Formula = '=IF(H11="","",TEXT(H11,"0.00")&" ("&TEXT(H11/$E11,"0.00")&")")'

$raw.Columns("B").ColumnWidth = 27.5

# ---------------------------------------------------------------------
# Sheet "collect": insert the matching summary row right after the
# collect_map_filter.rs row (row 5), before the collect_map.rs row
# (old row 6).
# ---------------------------------------------------------------------
$collect = $wb.Worksheets.Item("collect")
$collect.Rows("6:6").Insert()

$collect.Range("A6").Value2 = "drain_vec_collect_map_filter.rs"
$collect.Range("B6").Formula = '="[' + [char]0x21E8 + '](https://github.com/orxfun/orx-parallel/blob/main/benches/"&A6&")"'
$collect.Range("C6").Value2 = "``inputs.par_drain()`n  .map(map).filter(filter).collect()``"
$collect.Range("D6").Value2 = "58.37 (1.00)"
$collect.Range("E6").Value2 = "23.10 (0.40)"
$collect.Range("F6").Value2 = "13.96 (0.24)"
$collect.Range("G6").Value2 = "**12.84 (0.22)**"
$collect.Range("I6").Formula = '="**"&G6&"**"'
$collect.Rows(6).RowHeight = 28.8

$collect.Columns("A").ColumnWidth = 30

# ---------------------------------------------------------------------
# View state: the workbook was re-saved with "collect" as the active
# sheet (it used to be "early-exit"). Selections on Raw/early-exit are
# applied first so that activating "collect" last leaves it as the
# single tabSelected sheet.
# ---------------------------------------------------------------------
$raw.Range("I11:L11").Select()

$earlyExit = $wb.Worksheets.Item("early-exit")
$earlyExit.Range("E32").Select()

$collect.Activate()
$collect.Range("H21").Select()
